# Apply "New Metabuli version results" updates to the Error_types_table sheet.
# The sheet stores raw E (count), F (total), G (percentage = ROUND(E/F*100,2)) values
# for each Metabuli row; update them to reflect the new Metabuli run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-EF {
    param($Row, $E, $F)
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $g = [Math]::Round(($E / $F) * 100, 2)
    $ws.Cells.Item($Row, 7).Value = $g
}

function Set-DEF {
    param($Row, $D, $E, $F)
    $ws.Cells.Item($Row, 4).Value = $D
    Set-EF $Row $E $F
}

function Set-BDEF {
    param($Row, $B, $D, $E, $F)
    $ws.Cells.Item($Row, 2).Value = $B
    Set-DEF $Row $D $E $F
}

Set-EF 205 56 69
Set-EF 206 8  69
Set-EF 207 5  69
Set-EF 208 53 71
Set-EF 209 14 71
Set-EF 210 4  71
Set-EF 211 56 59
Set-EF 212 2  59
Set-EF 213 1  59
Set-EF 214 16 17
Set-EF 215 1  17
Set-EF 217 10 24
Set-EF 218 3  24
Set-EF 219 23 23
Set-BDEF 220 "Wadjemup" "Genus and species correct" 90 92
Set-DEF 221 "Genus correct, species wrong" 1 92
Set-DEF 222 "Genus and species wrong" 1 92
Set-EF 223 65 76
Set-EF 224 6  76
Set-EF 225 5  76
Set-EF 226 69 87
Set-EF 227 10 87
Set-EF 228 8  87
